# Week 15 simulations update
$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 15 play-by-play yardage data ---
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = "-1 2 5 2 9 2 2 -2 8 7 3 0 5 4 6 6 2 2 12 3 8 0 2 4 70 3 -1 2 0 5 9 3 2 7 0 1 -6 0 3 3 16 2 2 5 0 5 8 -1 3 14 2 0 3 6 5 7 3 8 4 11 0 1 14 2 2 0 13 4 13 1 10 0 -1 3 2 -3 10 4 7 2 -1 1 -2 1 0 -2 9 2 2 6 -1 7 8 1 4 14 3 2 31 3 3 3 3 -3 6 10 -1 12 1 5 2 -4 8 0 49 2 7 1 1 2 3 10 -1 11 4 5 10 1 9 3 11 3 4 0 3 0 5 8 4 0 30 4 2 4 4 1 1 0 6 2 1 6 0 2 7 3 3 5 2 0 7 2 8 4 2 6 2 2 7 1 1 7 2 -1 5 3 0 11 14 1 7 1 -2 -2 11 20 -3 7 13 1 -2 3 17 0 12 2 2 11 2 2 2 30 2 1 1 7 14 6 7 6 4 -4 2 0 1 10 3 2 0 2 1 -1 1 3 9 11 2 15 7 1 1 4 20 -4 8 9 1 2 12 2 10 2 4 2 4 2 -2 11 12 8 16 -1 9 1 2 2 6 6 2 2 0 3 0 7 4 3 8 4 4 3 9 -1 1 14 2 9 6 7 2 14 3 0 1 11 3 5 5 -1 5 4 8 9 7 2 6 2 16 1 11 5 3 1 13 9 1 5 6 3 3 9 4 1 3 2 2 1 3 3 -1 6 16 5 5 5 9 14 -1 1 1 6 8 3 4 -1 10 5 1"
$ydsWs.Range("B3").Value = "15 13 7 6 1 9 5 25 -4 5 11 6 2 6 17 16 2 10 6 8 20 7 22 13 3 6 13 14 8 6 6 7 4 10 16 33 15 12 9 55 14 5 10 11 9 0 4 11 14 6 6 10 8 5 28 6 17 6 9 27 6 21 -3 9 7 24 14 4 31 11 19 3 3 12 14 3 11 11 8 3 5 6 13 -3 9 9 32 9 5 10 9 6 14 12 12 3 12 19 10 9 15 9 39 23 6 9 5 23 6 7 14 23 5 1 3 6 11 7 5 7 10 14 3 14 11 6 10 5 7 26 13 12 12 12 26 6 18 16 6 4 4 1 6 6 8 31 1 0 6 6 9 12 13 4 11 9 8 8 5 2 9 10 3 7 9 10 1 31 5 31 12 19 19 15 7 3 4 -4 9 2 11 14 9 18 7 -3 25 5 8 11 44 19 7 6 9 17 6 3 11 20 8 -3 16 40 -2 7 5 64 4 13 15 26 10 7 4 10 6 15 7 5 -9 0 13 9 14 3 3 13 19 5 7 12 7 5 2 12 1 12 14 42 1 36 10 10 6 3 2 34 5 9 5 8 11 9 13 19 12 6 19 10 7 13 32 8 9 8 6 7 7 4 8 7 9 10 16 5 5 4"
$ydsWs.Range("C2").Value = "5 -2 4 2 -1 1 5 2 3 7 3 5 2 4 2 6 1 4 4 4 2 0 6 2 2 1 6 4 -1 12 7 11 5 4 1 4 2 2 2 13 8 5 4 2 -7 3 -1 7 5 2 1 1 3 4 5 7 11 2 1 -1 6 -2 3 4 3 5 4 0 2 3 5 3 5 5 4 5 6 3 3 3 5 7 6 5 2 4 5 20 -3 3 5 3 0 0 1 3 13 -2 5 18 3 8 -1 -1 4 9 2 1 4 4 1 4 6 2 6 4 3 9 3 6 18 6 2 0 1 3 1 5 3 1 3 0 20 10 4 12 1 9 2 1 2 5 10 11 2 3 3 2 0 9 1 6 1 9 5 -2 20 6 3 7 3 5 8 0 2 10 4 3 2 2 5 3 6 5 0 4 4 7 0 4 5 13 12 12 4 4 4 5 8 1 7 11 7 2 3 -2 7 2 5 10 9 10 17 -2 0 1 4 3 -8 7 4 2 5 11 31 0 3 4 3 0 5 -1 5 21 -2 15 5 25 1 7 5 3 1 2 6 1 23 1 1 7 5 6 3 3 5 1 5 -1 2 13 5 2 1 13 10 2 3 5 10 4 2 12 4 2 0 5 4 2 4 4 4 4 6 7 -4 1 5 1 13 35 5 9 3 4 4 2 3 5 4 -3 7 3 3 6 -3 8"
$ydsWs.Range("C3").Value = "42 -6 4 17 1 6 8 37 2 15 20 8 17 14 6 16 16 6 9 15 8 6 4 24 7 13 25 -2 -8 6 5 19 0 6 6 13 11 4 8 7 6 5 6 -5 14 5 3 22 7 4 10 23 4 13 13 20 6 12 24 12 9 49 32 7 13 14 4 8 13 7 8 17 4 7 14 32 4 23 50 4 5 13 6 59 5 14 14 3 18 1 22 16 9 48 2 3 25 11 16 6 8 8 10 33 31 29 3 40 8 51 7 34 4 5 14 8 4 6 9 15 6 12 -2 6 17 16 18 1 4 8 7 8 27 17 2 13 -2 20 4 12 1 17 14 7 5 14 4 23 20 19 5 22 14 3 1 5 0 32 17 21 4 18 5 4 24 9 32 5 10 7 4 20 11 4 16 18 3 4 6 24 36 8 6 4 5 9 11 12 14 2 7 24 -2 10 5 9 13 6 18 14 12 9 11 13 17 8 8 4 2 7 13 18 34 9 5 -2 16 9 9 12 13 9 29 14 -7 12 6 14 14 38 8 4 1 13 7 7 2 19 7 5 20 11 -2 9 9 11 13 4 6 14 8 14 7 2"

# --- OFF sheet: updated season totals ---
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = 155
$offWs.Range("D2").Value = 10
$offWs.Range("E2").Value = 16
$offWs.Range("F2").Value = 53
$offWs.Range("G2").Value = 50
$offWs.Range("J2").Value = 20
$offWs.Range("L2").Value = 222
$offWs.Range("M2").Value = 147
$offWs.Range("O2").Value = 22
$offWs.Range("P2").Value = 14
$offWs.Range("Q2").Value = 429
$offWs.Range("B3").Value = 8
$offWs.Range("C3").Value = 139
$offWs.Range("F3").Value = 74
$offWs.Range("G3").Value = 23
$offWs.Range("H3").Value = 29
$offWs.Range("I3").Value = 58
$offWs.Range("J3").Value = 33
$offWs.Range("N3").Value = 19

# --- DEF sheet: updated season totals ---
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 141
$defWs.Range("E2").Value = 9
$defWs.Range("F2").Value = 52
$defWs.Range("G2").Value = 42
$defWs.Range("I2").Value = 6
$defWs.Range("J2").Value = 26
$defWs.Range("L2").Value = 244
$defWs.Range("M2").Value = 151
$defWs.Range("O2").Value = 23
$defWs.Range("P2").Value = 8
$defWs.Range("Q2").Value = 430
$defWs.Range("C3").Value = 127
$defWs.Range("F3").Value = 65
$defWs.Range("G3").Value = 31
$defWs.Range("H3").Value = 25
$defWs.Range("I3").Value = 40
$defWs.Range("J3").Value = 52
$defWs.Range("N3").Value = 21

# --- ST sheet: updated special-teams totals + appended week-by-week arrays ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 62
$stWs.Range("D2").Value = 48
$stWs.Range("F2").Value = 253
$stWs.Range("G2").Value = 246
$stWs.Range("N2").Value = 54
$stWs.Range("O2").Value = 29
$stWs.Range("B3").Value = 52
$stWs.Range("D3").Value = "37 41 41 45 48 46 43 53 26 50 47 57 49 39 56 44 28 68 53 52 52 39 33 47 33 46 38 46 59 47 50 44 48 43 41 56 60 44 52 42 43 40 56 56 48 47 41 42"
$stWs.Range("D4").Value = "0 2 0 2 0 0 0 0 0 12 0 0 10 0 42 0 0 0 0 0 0 0 0 0 5 0 0 3 10 0 0 12 0 0 0 0 10 0 9 0 0 0 10 3 0 0 0 -4"
$stWs.Range("D5").Value = "0 0 2 11"

# --- TURNS sheet: updated totals ---
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("C2").Value = 8
$turnsWs.Range("E2").Value = 5
$turnsWs.Range("D3").Value = 6

# --- PEN sheet: updated totals ---
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B3").Value = 12

# --- Select ST sheet, matching the active tab / selection in the saved file ---
$stWs.Activate()
$stWs.Range("D6").Select() | Out-Null
